# Add the required "Experimental" boolean value (rendered as the literal
# text "true") to the ValueSet metadata sheet, and refresh the "Date"
# metadata value to reflect the new publish timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# B7 sits next to the "Experimental" label (A7) and was previously blank.
# Writing a bare "true" gets auto-typed as a Boolean by Excel, so we enter
# it with a leading apostrophe to force text, then strip the resulting
# quote-prefix formatting by re-pasting the original cell format from the
# neighboring "Experimental" label cell (A7) - this keeps the cell's style
# identical to before while leaving the stored value as the plain text
# string "true".
$ws.Range("B7").Value = "'true"
$ws.Range("A7").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# B8 holds the ISO-8601 "Date" metadata value; bump it to the new
# publication timestamp.
$ws.Range("B8").Value = "2023-02-01T09:05:11-06:00"
